$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Block A: rows 31..87 — every row's data shifts DOWN by one position
# (new row R = old row R-1), processed from the bottom (87) upward to
# the top (32) so that each source row is still in its original state
# when it is read. Row 31 itself is then overwritten with brand-new
# data (not derived from any existing row).
# ---------------------------------------------------------------------
for ($r = 87; $r -ge 32; $r--) {
    for ($c = 1; $c -le 18; $c++) {
        $srcCell = $ws.Cells.Item($r - 1, $c)
        $val = $srcCell.Value2
        $ws.Cells.Item($r, $c).Value2 = $val
    }
}

# New data for row 31
$ws.Cells.Item(31, 4).Value2  = 45061              # D31 Fecha
$ws.Cells.Item(31, 8).Value2  = "Cristal"            # H31 Variedad
$ws.Cells.Item(31, 11).Value2 = 16000                # K31 Precio minimo
$ws.Cells.Item(31, 12).Value2 = 16000                # L31 Precio maximo
$ws.Cells.Item(31, 13).Value2 = 16000                # M31 Precio promedio ponderado
$ws.Cells.Item(31, 14).Value2 = "$/saco 25 kilos"    # N31 Unidad de comercializacion
$ws.Cells.Item(31, 15).Value2 = "Región del Maule"  # O31 Origen
$ws.Cells.Item(31, 16).Value2 = 640                  # P31 Precio $/Kg

# ---------------------------------------------------------------------
# Block B: rows 180..206 — every row's data shifts UP by one position
# (new row R = old row R+1), processed from the top (180) downward to
# the bottom (206) so that each source row is still in its original
# state when it is read. Row 207 (now redundant, its data having moved
# into row 206) is then deleted outright.
# ---------------------------------------------------------------------
for ($r = 180; $r -le 206; $r++) {
    for ($c = 1; $c -le 18; $c++) {
        $srcCell = $ws.Cells.Item($r + 1, $c)
        $val = $srcCell.Value2
        $ws.Cells.Item($r, $c).Value2 = $val
    }
}

$ws.Rows.Item(207).Delete()
